$d = $word.ActiveDocument
$rng = $d.Content
$found = $rng.Find.Execute("vaisseau se déplace en consommant du carburant", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
Write-Host "found: $found"
Write-Host "rng text: $($rng.Text)"
$rng.HighlightColorIndex = 4
Write-Host "done"
